$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 2506
$ws1.Range("F3").Value = 537
$ws1.Range("F4").Value = 448
$ws1.Range("F5").Value = 276
$ws1.Range("F7").Value = 441
$ws1.Range("F8").Value = 1166
$ws1.Range("F9").Value = 523
$ws1.Range("F10").Value = 278
$ws1.Range("F12").Value = 334
$ws1.Range("F13").Value = 5384
$ws1.Range("F14").Value = 42
$ws1.Range("F15").Value = 1591
$ws1.Range("F16").Value = 3939
$ws1.Range("F17").Value = 383
$ws1.Range("F18").Value = 238
$ws1.Range("F20").Value = 4421
$ws1.Range("F21").Value = 5827
$ws1.Range("F22").Value = 137
$ws1.Range("F24").Value = 636
$ws1.Range("F25").Value = 3616
$ws1.Range("F26").Value = 450
$ws1.Range("F28").Value = 174
$ws1.Range("F29").Value = 114
$ws1.Range("F30").Value = 948
$ws1.Range("F31").Value = 1322
$ws1.Range("F32").Value = 141
$ws1.Range("F33").Value = 191
$ws1.Range("F34").Value = 1536
$ws1.Range("F35").Value = 178
$ws1.Range("F36").Value = 1600
$ws1.Range("F37").Value = 141
$ws1.Range("F38").Value = 1047
$ws1.Range("F40").Value = 1344
$ws1.Range("F41").Value = 590
$ws1.Range("F42").Value = 81
$ws1.Range("F43").Value = 170
$ws1.Range("F44").Value = 2704
$ws1.Range("F45").Value = 116
$ws1.Range("F46").Value = 242
$ws1.Range("F47").Value = 393
$ws1.Range("F49").Value = 3838

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 4
$ws2.Range("F5").Value = 1145
$ws2.Range("F12").Value = 16

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 3593

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 3593
$ws4.Range("F3").Value = 2506
$ws4.Range("F4").Value = 537
$ws4.Range("F5").Value = 448
$ws4.Range("F6").Value = 276
$ws4.Range("F7").Value = 1145
$ws4.Range("F9").Value = 441
$ws4.Range("F10").Value = 1166
$ws4.Range("F11").Value = 523
$ws4.Range("F12").Value = 278
$ws4.Range("F14").Value = 334
$ws4.Range("F15").Value = 5384
$ws4.Range("F16").Value = 16
$ws4.Range("F17").Value = 1591
$ws4.Range("F18").Value = 4421
$ws4.Range("F19").Value = 5827
$ws4.Range("F20").Value = 137
$ws4.Range("F22").Value = 636
$ws4.Range("F23").Value = 3616
$ws4.Range("F24").Value = 450
$ws4.Range("F26").Value = 174
$ws4.Range("F27").Value = 114
$ws4.Range("F28").Value = 948
$ws4.Range("F29").Value = 1323
$ws4.Range("F30").Value = 142
$ws4.Range("F31").Value = 192
$ws4.Range("F32").Value = 1536
$ws4.Range("F33").Value = 178
$ws4.Range("F34").Value = 1600
$ws4.Range("F36").Value = 1047
$ws4.Range("F38").Value = 590
$ws4.Range("F41").Value = 81
$ws4.Range("F43").Value = 2704
$ws4.Range("F45").Value = 116
$ws4.Range("F46").Value = 242
$ws4.Range("F47").Value = 393
$ws4.Range("F49").Value = 3838
